$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 2.88
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
$ws.Range("P3").Value = 1.53
$ws.Range("Q3").Value = 2.38
$ws.Range("Z3").Value = 7.5
$ws.Range("AA3").Value = 6
$ws.Range("AD3").Value = 451

# Row 5
$ws.Range("G5").Value = 1.62

# Row 9
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 3.6
$ws.Range("T9").Value = 8
$ws.Range("U9").Value = 9
$ws.Range("X9").Value = 15
$ws.Range("AB9").Value = 15
$ws.Range("AD9").Value = 201
$ws.Range("AE9").Value = 12

# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("N12").Value = 2.03
$ws.Range("O12").Value = 1.83
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 11
$ws.Range("AA12").Value = 6.5
$ws.Range("AB12").Value = 15

# Row 14
$ws.Range("G14").Value = 1.62
$ws.Range("H14").Value = 4.1
$ws.Range("I14").Value = 5
$ws.Range("P14").Value = 1.33
$ws.Range("Q14").Value = 3.25
$ws.Range("R14").Value = 1.7
$ws.Range("S14").Value = 2.05
$ws.Range("U14").Value = 8.5
$ws.Range("X14").Value = 12
$ws.Range("AA14").Value = 8
$ws.Range("AB14").Value = 15
$ws.Range("AF14").Value = 29
$ws.Range("AI14").Value = 41
$ws.Range("AJ14").Value = 41

# Row 16
$ws.Range("N16").Value = 1.25
$ws.Range("O16").Value = 3.55
$ws.Range("R16").Value = 2.41
$ws.Range("S16").Value = 1.5

# Row 17
$ws.Range("K17").Value = 9
$ws.Range("P17").Value = 1.41
$ws.Range("Q17").Value = 2.62

# Row 18
$ws.Range("P18").Value = 1.22

# Row 19
$ws.Range("L19").Value = 1.24
$ws.Range("Q19").Value = 2.92
$ws.Range("S19").Value = 2.1

# Row 20
$ws.Range("K20").Value = 10
$ws.Range("P20").Value = 1.37

# Row 23
$ws.Range("N23").Value = 2.1
$ws.Range("O23").Value = 1.7
